$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 2 (new accelerometer samples
# captured earlier than the existing series) and clear any inherited
# formatting so the new rows look like plain data rows.
$ws.Rows.Item(2).Resize(3).Insert()
$ws.Range("A2:C4").ClearFormats()

$topRows = @(
    @(-0.084752082824707, 0.6204710006713867, -1.197814345359802),
    @(-0.5587072372436523, 0.5920883417129517, -0.9495211839675904),
    @(-0.2052898406982422, 0.80674147605896, -1.046440482139587)
)

for ($i = 0; $i -lt $topRows.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 1).Value2 = $topRows[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $topRows[$i][1]
    $ws.Cells.Item($r, 3).Value2 = $topRows[$i][2]
}

# Append 7 more rows of new samples after the existing data (which now
# ends at row 24).
$bottomRows = @(
    @(0.0039987564086914, 0.546174168586731, -0.7374091148376465),
    @(0.0831842422485351, 0.5668889284133911, -0.8130950927734375),
    @(-0.0142126083374023, 0.51572585105896, -0.7260744571685791),
    @(0.1546173095703125, 0.5381616353988647, -0.7814648151397705),
    @(0.2052326202392578, 0.5754936933517456, -0.8383152484893799),
    @(-0.07891082763671869, 0.5558477640151978, -0.7180624008178711),
    @(0.1948976516723632, 0.6977589726448059, -0.9572491645812988)
)

for ($i = 0; $i -lt $bottomRows.Count; $i++) {
    $r = 25 + $i
    $ws.Cells.Item($r, 1).Value2 = $bottomRows[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $bottomRows[$i][1]
    $ws.Cells.Item($r, 3).Value2 = $bottomRows[$i][2]
}

"done"
